$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Type the values in the order that reproduces the original shared-strings
# table layout (MD5, hash1, hash2, hash3, hash4), then rearrange the cells
# into their final order. Re-assigning a value that is already present in
# the shared-string table does not create a new entry, so the final table
# order matches what a user would get by typing the labels once and then
# sorting/reordering the rows.
$ws.Range("A1").Value = "MD5"
$ws.Range("A2").Value = "5AF1409B22F39CDFDA5906D61B6266A9"
$ws.Range("A3").Value = "7D54F8AE69FE019A7D744338DE97885B"
$ws.Range("A4").Value = "5F4DCC3B5AA765D61D8327DEB882CF99"
$ws.Range("A5").Value = "69904329D6E12BFAF68602EAA1E4EEB1"

$ws.Range("A2").Value = "5F4DCC3B5AA765D61D8327DEB882CF99"
$ws.Range("A3").Value = "7D54F8AE69FE019A7D744338DE97885B"
$ws.Range("A4").Value = "69904329D6E12BFAF68602EAA1E4EEB1"
$ws.Range("A5").Value = "5AF1409B22F39CDFDA5906D61B6266A9"

# Size column A to fit the longest value (the 32-character MD5 strings).
$ws.Columns(1).AutoFit() | Out-Null

# Match the saved print-header/footer margins.
$ws.PageSetup.HeaderMargin = 0.51166665554046631 * 72
$ws.PageSetup.FooterMargin = 0.51166665554046631 * 72

# Select the full populated range, matching the saved selection.
$ws.Range("A1:A5").Select() | Out-Null
